$wb = $excel.ActiveWorkbook

# --- entities sheet: delete row 31 (entire row), shifting subsequent rows up ---
$wsEntities = $wb.Worksheets.Item("entities")
$null = $wsEntities.Rows.Item(31).Delete()

# --- identifiers sheet: move selection to B3 ---
$wsIdentifiers = $wb.Worksheets.Item("identifiers")
$null = $wsIdentifiers.Range("B3").Select()

# --- iterations sheet: move selection to E6 ---
$wsIterations = $wb.Worksheets.Item("iterations")
$null = $wsIterations.Range("E6").Select()

# --- make "entities" the active/selected tab (was "district") ---
$null = $wsEntities.Activate()
$null = $wsEntities.Range("K1").Select()
